# Insert a brand-new introductory paragraph describing the Random Forest
# model, ahead of the existing first paragraph ("The first iteration of
# training the Random Forest Classifier ...").

$d = $word.ActiveDocument

# Anchor on the existing first paragraph and push a new empty paragraph
# in front of it -- this leaves the original paragraph's own
# formatting/rsids completely untouched.
$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphBefore()

# The freshly-minted paragraph is now Paragraphs(1); build up its text
# piece by piece, mirroring the distinct runs the sentence was actually
# composed from.
$newPara = $d.Paragraphs.First
$ins = $newPara.Range
$ins.Collapse(1)

$runs = @(
    "Random Forest is an ‘ensemble’ model that fits based on majority voting from numerous decision trees",
    " which corrects for overfitting",
    ". ",
    "As Random Forest is a non-parametric algorithm, it ",
    "requires little data preparation beforehand.",
    " ",
    "V",
    "ariables can be ranked according to importance based on Gini index, though ",
    "how a variable affects final output ",
    "is less interpretable than logistic regression.",
    " "
)

# Insert each chunk of text as its own run. A transient bookmark is
# dropped between successive chunks so the engine treats them as
# separate insertions instead of silently re-coalescing adjacent,
# identically-formatted runs; the bookmarks are removed again once all
# the text is in place, leaving no trace in the final document.
$markNames = New-Object System.Collections.ArrayList
for ($i = 0; $i -lt $runs.Length; $i++) {
    $ins.InsertAfter($runs[$i])
    $ins.Collapse(0)
    if ($i -lt $runs.Length - 1) {
        $markName = "zzTmpRunSplit" + $i
        $d.Bookmarks.Add($markName, $ins) | Out-Null
        $markNames.Add($markName) | Out-Null
    }
}
foreach ($markName in $markNames) {
    $d.Bookmarks($markName).Delete()
}

Write-Output "inserted intro paragraph"
